$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "24.384.31"
$ws.Range("E2").Value = "  +8.75%  "
Set-TextCell $ws.Range("D3") "1.678.95"
$ws.Range("E3").Value = "  +4.06%  "
Set-TextCell $ws.Range("D4") "1.008"
$ws.Range("E4").Value = "  +0.30%  "
Set-TextCell $ws.Range("D5") "307.34"
$ws.Range("E5").Value = "  +0.86%  "
Set-TextCell $ws.Range("D6") "0.9993"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  +0.50%  "
Set-TextCell $ws.Range("D8") "0.3439"
$ws.Range("E8").Value = "  +0.59%  "
Set-TextCell $ws.Range("D9") "47.72"
$ws.Range("E9").Value = "  +12.81%  "
Set-TextCell $ws.Range("D10") "1.170"
$ws.Range("E10").Value = "  +2.84%  "
Set-TextCell $ws.Range("D11") "0.07274"
$ws.Range("E11").Value = "  +2.64%  "
Set-TextCell $ws.Range("D12") "1.000"
$ws.Range("E12").Value = "  -0.18%  "
Set-TextCell $ws.Range("D13") "6.117"
$ws.Range("E13").Value = "  +2.92%  "
Set-TextCell $ws.Range("D14") "20.27"
$ws.Range("E14").Value = "  +2.40%  "
Set-TextCell $ws.Range("D15") "6.751"
$ws.Range("E15").Value = "  +1.29%  "
Set-TextCell $ws.Range("D16") "1.681.40"
$ws.Range("E16").Value = "  +4.65%  "
Set-TextCell $ws.Range("D17") "0.00001107"
$ws.Range("E17").Value = "  +1.67%  "
Set-TextCell $ws.Range("D18") "0.9988"
$ws.Range("E18").Value = "  +0.74%  "
Set-TextCell $ws.Range("D19") "0.06681"
$ws.Range("E19").Value = "  -1.48%  "
Set-TextCell $ws.Range("D20") "81.03"
$ws.Range("E20").Value = "  +3.32%  "
Set-TextCell $ws.Range("D21") "16.44"
$ws.Range("E21").Value = "  +2.10%  "
Set-TextCell $ws.Range("D22") "6.114"
$ws.Range("E22").Value = "  +0.96%  "
Set-TextCell $ws.Range("D23") "12.19"
$ws.Range("E23").Value = "  +2.47%  "
Set-TextCell $ws.Range("D24") "24.376.91"
$ws.Range("E24").Value = "  +8.38%  "
Set-TextCell $ws.Range("D25") "2.455"
$ws.Range("E25").Value = "  +2.77%  "
Set-TextCell $ws.Range("D26") "2.660"
$ws.Range("E26").Value = "  +4.36%  "
Set-TextCell $ws.Range("D27") "153.65"
$ws.Range("E27").Value = "  +2.30%  "
Set-TextCell $ws.Range("D28") "19.52"
Set-TextCell $ws.Range("D29") "1.867.23"
$ws.Range("E29").Value = "  +4.49%  "
Set-TextCell $ws.Range("D30") "127.38"
$ws.Range("E30").Value = "  +3.30%  "
Set-TextCell $ws.Range("D31") "6.314"
$ws.Range("E31").Value = "  +2.73%  "
Set-TextCell $ws.Range("D32") "4.056"
$ws.Range("E32").Value = "  +0.11%  "
Set-TextCell $ws.Range("D33") "0.9749"
$ws.Range("E33").Value = "  +1.88%  "
Set-TextCell $ws.Range("D34") "0.08476"
$ws.Range("E34").Value = "  +2.63%  "
Set-TextCell $ws.Range("D35") "1.706"
$ws.Range("E35").Value = "  +3.16%  "
Set-TextCell $ws.Range("D36") "12.37"
$ws.Range("E36").Value = "  +3.01%  "
Set-TextCell $ws.Range("D37") "0.06497"
$ws.Range("E37").Value = "  +6.36%  "
Set-TextCell $ws.Range("D38") "5.358"
$ws.Range("E38").Value = "  +1.56%  "
Set-TextCell $ws.Range("D39") "8.875"
$ws.Range("E39").Value = "  +2.95%  "
Set-TextCell $ws.Range("D40") "0.02326"
$ws.Range("E40").Value = "  +4.16%  "
Set-TextCell $ws.Range("D41") "1.256"
$ws.Range("E41").Value = "  -1.20%  "
Set-TextCell $ws.Range("D42") "0.2103"
$ws.Range("E42").Value = "  +3.72%  "
Set-TextCell $ws.Range("D43") "0.6157"
$ws.Range("E43").Value = "  +3.63%  "
Set-TextCell $ws.Range("D44") "0.9981"
$ws.Range("E44").Value = "  +0.71%  "
Set-TextCell $ws.Range("D47") "0.5940"
$ws.Range("E47").Value = "  +3.86%  "
Set-TextCell $ws.Range("D48") "127.39"
$ws.Range("E48").Value = "  +0.22%  "
Set-TextCell $ws.Range("D49") "2.024"
$ws.Range("E49").Value = "  +1.93%  "
Set-TextCell $ws.Range("D50") "0.07190"
$ws.Range("E50").Value = "  +5.41%  "
Set-TextCell $ws.Range("D51") "75.93"
$ws.Range("E51").Value = "  +2.56%  "

# Row 45/46 swap: PancakeSwap <-> EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D45") "13.23"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Range("D46") "3.775"
$ws.Range("E46").Value = "  -1.27%  "
